# Scheduled runner update: refresh market price/profit columns (H-N)
# across the Leve worksheets per latest snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 206
$ws.Range("I33").Value = 120
$ws.Range("J33").Value = 292
$ws.Range("K33").Value = 120
$ws.Range("L33").Value = 292
$ws.Range("M33").Value = 109
$ws.Range("N33").Value = -750

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1441.7222
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1441.7222
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1441.7222
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -1579.7222

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3217.5881
$ws.Range("I76").Value = 3212.4375
$ws.Range("K76").Value = 3212.4375
$ws.Range("M76").Value = -2897.4375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3217.5881
$ws.Range("I79").Value = 3212.4375
$ws.Range("K79").Value = 3212.4375
$ws.Range("M79").Value = -2120.4375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3856468
$ws.Range("I137").Value = 14287229
$ws.Range("J137").Value = 1621.4565
$ws.Range("K137").Value = 42861687
$ws.Range("L137").Value = 4864.3695
$ws.Range("M137").Value = -42859137
$ws.Range("N137").Value = -9964.369500000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3300.96
$ws.Range("I138").Value = 1304.6
$ws.Range("J138").Value = 4026.9092
$ws.Range("K138").Value = 3913.8
$ws.Range("L138").Value = 12080.7276
$ws.Range("M138").Value = 1226.2
$ws.Range("N138").Value = -22360.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 7556.4
$ws.Range("I31").Value = 7556.4
$ws.Range("K31").Value = 7556.4
$ws.Range("M31").Value = -7262.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 22225388
$ws.Range("I122").Value = 3955.75
$ws.Range("K122").Value = 11867.25
$ws.Range("M122").Value = -9417.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 136017.73
$ws.Range("I132").Value = 145251.42
$ws.Range("J132").Value = 127938.25
$ws.Range("K132").Value = 435754.26
$ws.Range("L132").Value = 383814.75
$ws.Range("M132").Value = -433224.26
$ws.Range("N132").Value = -388874.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 7952
$ws.Range("J81").Value = 7952
$ws.Range("L81").Value = 7952
$ws.Range("N81").Value = -10074

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 7952
$ws.Range("J84").Value = 7952
$ws.Range("L84").Value = 23856
$ws.Range("N84").Value = -34464

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 834.40625
$ws.Range("I105").Value = 888
$ws.Range("J105").Value = 643
$ws.Range("K105").Value = 888
$ws.Range("L105").Value = 643
$ws.Range("M105").Value = 859
$ws.Range("N105").Value = -4137

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1197.2222
$ws.Range("I107").Value = 493.75
$ws.Range("J107").Value = 1760
$ws.Range("K107").Value = 493.75
$ws.Range("L107").Value = 1760
$ws.Range("M107").Value = 1426.25
$ws.Range("N107").Value = -5600

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 22923.895
$ws.Range("I132").Value = 1201.1562
$ws.Range("J132").Value = 69265.734
$ws.Range("K132").Value = 3603.4686
$ws.Range("L132").Value = 207797.202
$ws.Range("M132").Value = -1073.4686
$ws.Range("N132").Value = -212857.202

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 969.2
$ws.Range("I68").Value = 495.33334
$ws.Range("J68").Value = 1033.8182
$ws.Range("K68").Value = 1486.00002
$ws.Range("L68").Value = 3101.4546
$ws.Range("M68").Value = -675.0000199999999
$ws.Range("N68").Value = -4723.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 969.2
$ws.Range("I71").Value = 495.33334
$ws.Range("J71").Value = 1033.8182
$ws.Range("K71").Value = 4458.00006
$ws.Range("L71").Value = 9304.363799999999
$ws.Range("M71").Value = -402.0000600000003
$ws.Range("N71").Value = -17416.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50507.184
$ws.Range("I70").Value = 71453.266
$ws.Range("J70").Value = 5622.7144
$ws.Range("K70").Value = 71453.266
$ws.Range("L70").Value = 5622.7144
$ws.Range("M70").Value = -71183.266
$ws.Range("N70").Value = -6162.7144

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 50507.184
$ws.Range("I73").Value = 71453.266
$ws.Range("J73").Value = 5622.7144
$ws.Range("K73").Value = 71453.266
$ws.Range("L73").Value = 5622.7144
$ws.Range("M73").Value = -70517.266
$ws.Range("N73").Value = -7494.7144

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3331.8333
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 3712.3572
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 3712.3572
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -5708.3572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3331.8333
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 3712.3572
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 18561.786
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -28545.786

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 862.6087
$ws.Range("I107").Value = 780.2308
$ws.Range("J107").Value = 969.7
$ws.Range("K107").Value = 780.2308
$ws.Range("L107").Value = 969.7
$ws.Range("M107").Value = 1139.7692
$ws.Range("N107").Value = -4809.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 65335.97
$ws.Range("I132").Value = 65501.75
$ws.Range("J132").Value = 65170.188
$ws.Range("K132").Value = 196505.25
$ws.Range("L132").Value = 195510.564
$ws.Range("M132").Value = -193975.25
$ws.Range("N132").Value = -200570.564

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2559.8
$ws.Range("I7").Value = 2574.75
$ws.Range("J7").Value = 2500
$ws.Range("K7").Value = 2574.75
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = -2462.75
$ws.Range("N7").Value = -2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3734
$ws.Range("I122").Value = 3734
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11202
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8752
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2559.8
$ws.Range("I126").Value = 2574.75
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 7724.25
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -5254.25
$ws.Range("N126").Value = -12440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 41185.27
$ws.Range("I132").Value = 2027.625
$ws.Range("J132").Value = 103837.5
$ws.Range("K132").Value = 6082.875
$ws.Range("L132").Value = 311512.5
$ws.Range("M132").Value = -3552.875
$ws.Range("N132").Value = -316572.5
